# This script applies the commit's formatting change:
#   - Several bold, centered heading paragraphs ("Условие задачи",
#     "Решении задачи", "Состав данных", the "Форма вывода" paragraph and
#     the "Листинг программы, ..." paragraph) get an explicit font size of
#     14pt (w:sz/w:szCs = 28 half-points) applied to every run (and to the
#     paragraph mark's run properties).
#   - A new, empty, centered/bold paragraph is inserted right after the
#     "Условие задачи" and "Решении задачи" headings (but not after the
#     other headings).

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $searchText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    $targetStart = $rng.Start
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $targetStart -and $p.Range.End -gt $targetStart) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphHeadingSize($doc, $searchText) {
    $idx = Find-ParagraphIndex $doc $searchText
    if ($idx -eq -1) {
        Write-Host "NOT FOUND:" $searchText
        return $null
    }
    $p = $doc.Paragraphs.Item($idx)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
    return $p
}

# 1) "Условие задачи" heading - first add a blank heading-style paragraph
#    right after it (while the run still has the original, smaller size),
#    then bump the heading's own size.
$d.Content.Find.Execute("Условие задачи", $true, $false, $false, $false, $false, $true, 1, $false, "Условие задачи^p", 2) | Out-Null
Set-ParagraphHeadingSize $d "Условие задачи" | Out-Null

# 2) "Решении задачи" heading - same treatment.
$d.Content.Find.Execute("Решении задачи", $true, $false, $false, $false, $false, $true, 1, $false, "Решении задачи^p", 2) | Out-Null
Set-ParagraphHeadingSize $d "Решении задачи" | Out-Null

# 3) "Состав данных" heading - bump size only (no new paragraph).
Set-ParagraphHeadingSize $d "Состав данных" | Out-Null

# 4) The "Форма вывода" paragraph (also holds the preceding drawing run) -
#    bump size across the whole paragraph (drawing run, "Форма ", "вывода"
#    and the trailing line break run).
Set-ParagraphHeadingSize $d "Форма" | Out-Null

# 5) The "Листинг программы, ..." paragraph - bump size across every run.
Set-ParagraphHeadingSize $d "Листинг программы" | Out-Null
